$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = '腹部逆位'
$ws.Range("G3").Value = '腹部閉鎖位'
$ws.Range("G4").Value = '左房異型症'
$ws.Range("G5").Value = '右房異数性'
$ws.Range("G6").Value = '気管支異数性'
$ws.Range("G7").Value = '過剰成長'
$ws.Range("G8").Value = '下肢半低身長症'
$ws.Range("G9").Value = '上肢半身異栄養症'
$ws.Range("G10").Value = '下肢半身肥大'
$ws.Range("G11").Value = '下肢半側肥大症'
$ws.Range("G12").Value = '上肢半身肥大'
$ws.Range("G13").Value = 'ウエスト・ヒップ比増加'
$ws.Range("G14").Value = 'ウエスト／ヒップ比低下'
$ws.Range("G15").Value = '妊娠月齢に比して大きい'
$ws.Range("G16").Value = '過体重'
$ws.Range("G17").Value = 'II度肥満'
$ws.Range("G18").Value = 'III度肥満'
$ws.Range("G19").Value = '腹部肥満'
$ws.Range("G20").Value = 'I度肥満'
$ws.Range("G21").Value = '小児期発症の三頭筋肥満'
$ws.Range("G22").Value = '細身体型'
$ws.Range("G23").Value = '反復性感染症による二次性発育不全'
$ws.Range("G24").Value = '重症成長障害'
$ws.Range("G25").Value = '乳児期の発育不全'
$ws.Range("G26").Value = '肥満度の低下'
$ws.Range("G27").Value = '妊娠月齢に比して小さい'
$ws.Range("G29").Value = '悪液質'
$ws.Range("G30").Value = '肥満度の増加'
$ws.Range("G31").Value = '上腕三頭筋皮下脂肪厚の増加'
$ws.Range("G32").Value = '非対称性低身長'
$ws.Range("G33").Value = '新生児短躯性低身長'
$ws.Range("G34").Value = '乳児期発症型体幹短小'
$ws.Range("G35").Value = '致死性体幹短小'
$ws.Range("G36").Value = '小児期発症短躯性低身長症'
$ws.Range("G37").Value = '致死性四肢短縮型低身長症'
$ws.Range("G38").Value = 'メソメリック性低身長症'
$ws.Range("G39").Value = '新生児短肢性低身長症'
$ws.Range("G40").Value = 'リゾ-メソ-アクロメリック四肢短縮症'
$ws.Range("G41").Value = 'リゾメリック腕短縮症'
$ws.Range("G42").Value = 'メソメリック/リゾメリック四肢短縮症'
$ws.Range("G43").Value = 'リゾメリック下肢短縮症'
$ws.Range("G44").Value = '重症短下肢小人症'
$ws.Range("G45").Value = '小児期発症短下肢小人症'
$ws.Range("G46").Value = '重症低身長症'
$ws.Range("G47").Value = '軽度低身長'
$ws.Range("G48").Value = '中等度低身長'
$ws.Range("G49").Value = '下垂体性小人症'
$ws.Range("G50").Value = '出生時の体長が3パーセンタイル未満'
$ws.Range("G51").Value = '過成長'
$ws.Range("G52").Value = '比例高身長'
$ws.Range("G53").Value = '出生時の体長が97パーセンタイル以上'
$ws.Range("G54").Value = '不釣り合いな高身長'
$ws.Range("G55").Value = '上位下位比の増加'
$ws.Range("G56").Value = '上下肢節比の減少'
$ws.Range("G57").Value = '成長ホルモンなしでの成長'
$ws.Range("G58").Value = '体脂肪率の増加'
$ws.Range("G60").Value = '初潮の遅れ'
$ws.Range("G61").Value = '副初期の遅延'
$ws.Range("G62").Value = '中等度の子宮内発育遅延'
$ws.Range("G63").Value = '軽度の子宮内発育遅延'
$ws.Range("G64").Value = '重度の子宮内発育遅延'
$ws.Range("G65").Value = '軽度の出生後発育遅延'
$ws.Range("G66").Value = '中等度の出生後発育遅延'
$ws.Range("G67").Value = '重度の出生後発育遅延'
$ws.Range("G68").Value = '思春期成長スパートの欠如'
$ws.Range("G69").Value = '体脂肪率の低下'

$null = $ws.Range("A1").EntireRow.Select()
